# "Slides 7 & 8 Officer Practice"
#
# The second slide in the deck (id 263) is a half-finished placeholder
# slide ("11-2" slide-number placeholder + Pearson copyright footer, plus
# an empty Title/Content/Content placeholder trio) that never got any
# real content. It gets removed entirely, so the deck goes from 4 slides
# down to 3: the blank opening slide, "The Last Bear Market", and
# "What do I do in a Bear Market".

$p = $ppt.ActivePresentation

# Remove the stray/unused slide (was slide index 2).
$p.Slides.Item(2).Delete()

# After the delete, the two remaining content slides shift up to
# positions 2 and 3. Give their Title / Content placeholders their
# normal default PowerPoint names (they were unnamed before).
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).Name = "Title 1"
$s2.Shapes.Item(2).Name = "Content Placeholder 2"

$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).Name = "Title 1"
$s3.Shapes.Item(2).Name = "Content Placeholder 2"
